$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)      # "2021-Q3"
$totalSheet = $wb.Worksheets.Item(2)  # "总计"

# --- Create the new "2022-Q1" sheet by duplicating the "总计" sheet so it
# inherits the same sheetPr / pageSetup / sheetFormatPr, then place it
# right after "2021-Q3" and rename it. (Re-fetch worksheet handles by name
# afterwards - this engine's worksheet object references can end up
# pointing at the wrong physical sheet if reused across a Copy() call.)
$totalSheet.Copy($null, $sheet1)
$q1Sheet = $wb.Worksheets.Item("总计 (2)")
$q1Sheet.Name = "2022-Q1"
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# Clear the copied "总计" sample data out of the new sheet before writing
# the real 2022-Q1 holdings.
$q1Sheet.Range("A1:D2").ClearContents

# Extend the header style (already present on B1:D1 from the copied sheet)
# across the new E1:H1 header cells before filling in their text, so we
# don't introduce any new style definitions.
$q1Sheet.Range("B1").Copy()
$q1Sheet.Range("E1:H1").PasteSpecial(-4122)

# --- Headers for "2022-Q1" ---
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# --- Data row for "2022-Q1" ---
$q1Sheet.Range("A2").Value = 0

# Text-like columns must stay text (not auto-converted to numbers), so
# force a text number format before writing, then strip the style back off
# again so the cells end up unstyled like the source file.
$q1Sheet.Range("B2:G2").NumberFormat = "@"
$q1Sheet.Range("B2").Value = "165524"
$q1Sheet.Range("C2").Value = "信诚中证智能家居指数（LOF）"
$q1Sheet.Range("D2").Value = "0.40"
$q1Sheet.Range("E2").Value = "93.89"
$q1Sheet.Range("F2").Value = "1.15"
$q1Sheet.Range("G2").Value = "0.0046"
$q1Sheet.Range("B2:G2").Style = "Normal"

# H2 is a real number.
$q1Sheet.Range("H2").Value = 10

# --- Update the "总计" sheet: insert a new row for 2022-Q1 above the
# existing 2021-Q3 row, shifting it down. Write the shifted row first so
# the style-bearing A column copies cleanly.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A2").Copy($totalSheet.Range("A3"))
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q3"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 0.04

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0
